$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert a new row at position 24 (shifts old rows 24-29 down to 25-30),
#    making room for a brand-new period entry (2508) at the bottom of the
#    "estado de cuenta" table, right before the signature block.
# ------------------------------------------------------------------
$ws.Rows.Item(24).Insert()

# Copy the special "closing" row formatting (thicker bottom border) that
# used to belong to row 23 down onto the new row 24 - it is now the last
# detail row of the table.
$ws.Range("B23:J23").Copy()
$ws.Range("B24:J24").PasteSpecial(-4122)

# Row 23 is no longer the last row of the table, so it goes back to using
# the regular (non-closing) row formatting, like row 22.
$ws.Range("B22:J22").Copy()
$ws.Range("B23:J23").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 2) Update the account summary figures.
# ------------------------------------------------------------------
$ws.Range("E11").Value = 472267   # VALOR MORA total
$ws.Range("F13").Value = 9        # Cant. Periodos

# ------------------------------------------------------------------
# 3) Re-date every period row: the table now starts with the oldest
#    period (2412) and lists each following period in ascending order,
#    ending with the newest period (2508) on the new closing row.
# ------------------------------------------------------------------
$ws.Range("E16").Value = "2412"
$ws.Range("F16").Value = 24267

$ws.Range("E17").Value = "2501"
$ws.Range("F17").Value = 56000

$ws.Range("E18").Value = "2502"
$ws.Range("F18").Value = 56000

$ws.Range("E19").Value = "2503"
$ws.Range("F19").Value = 56000

$ws.Range("E20").Value = "2504"
$ws.Range("F20").Value = 56000

$ws.Range("E21").Value = "2505"
$ws.Range("F21").Value = 56000

$ws.Range("E22").Value = "2506"
$ws.Range("F22").Value = 56000

$ws.Range("E23").Value = "2507"
$ws.Range("F23").Value = 56000
$ws.Range("G23").Value = 1400000

# ------------------------------------------------------------------
# 4) Fill in the brand-new period row (2508).
# ------------------------------------------------------------------
$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1007120425"
$ws.Range("D24").Value = "ABEL CARDALES MATOS"
$ws.Range("E24").Value = "2508"
$ws.Range("F24").Value = 56000
$ws.Range("G24").Value = 1400000
